$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: clear header text, keep the cells present (empty) with default style
$ws.Range("A1:G1").ClearContents()
$ws.Range("A1:G1").Style = "Normal"

# Row 2: "Входные данные" / variable entry
$ws.Range("A2").Value() = "Входные данные"
$ws.Range("B2").Value() = "variable"
$ws.Range("C2").Value() = "integer"
$ws.Range("D2").Value() = "Простая переменная"
$ws.Range("E2").Value() = "0-100"
$ws.Range("F2").Value() = "No format"
$ws.Range("G2").Value() = "does sth"

# Row 3: "Входные данные" / another_variable entry (new row replacing old HelloWorld row)
$ws.Range("A3").Value() = "Входные данные"
$ws.Range("B3").Value() = "another_variable"
$ws.Range("C3").Value() = "integer"
$ws.Range("D3").Value() = "Простая переменная"
$ws.Range("E3").Value() = "0-100"
$ws.Range("F3").Value() = "Format"
$ws.Range("G3").Value() = "a really long meaning"
